$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B18: was a number (1517492), becomes text "1517492"
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "1517492"

# New row 19 data
$ws.Range("A19").Value = 123444
$ws.Range("B19").Value = "PROPRO"
$ws.Range("C19").Value = "uuuuuuuuuuu"
$ws.Range("D19").Value = "Mètre"
$ws.Range("E19").Value = "Barre de 6m"
$ws.Range("F19").Value = 5
$ws.Range("G19").Value = 50
$ws.Range("H19").Value = "Site principal"
$ws.Range("I19").Value = "Stockage"
$ws.Range("J19").Value = "E3"
$ws.Range("K19").Value = "FournX"
$ws.Range("L19").Value = 45
$ws.Range("M19").Value = "Profilés"
$ws.Range("N19").Value = "Structure"
$ws.Range("O19").NumberFormat = "@"
$ws.Range("O19").Value = "2481023879"
$ws.Range("P19").Value = 40
$ws.Range("Q19").NumberFormat = "@"
$ws.Range("Q19").Value = "2025-05-28"

Write-Output "done"
